$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The workbook lists one account per row (Conta, Nome, Saldo) starting at row 2.
# Two new accounts need to be inserted right above account 004313254 (GUSTAVO),
# which currently sits on row 5 - so insert two blank rows there first, which
# pushes that row (and everything below it) down by two rows.
$insertRange = $ws.Range("A5:A6").EntireRow
$insertRange.Insert()

# "Conta" values carry significant leading zeros, so force these two cells to
# Text format before assigning them - otherwise Excel would interpret the
# digit string as a number and drop the leading zeros.
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).NumberFormat = "@"

# New row 1: account 004222784 / RAFAEL / 13000
$ws.Cells.Item(5, 1).Value = "004222784"
$ws.Cells.Item(5, 2).Value = "RAFAEL"
$ws.Cells.Item(5, 3).Value = 13000

# New row 2: account 004243043 / SUELI / 12752.5
$ws.Cells.Item(6, 1).Value = "004243043"
$ws.Cells.Item(6, 2).Value = "SUELI"
$ws.Cells.Item(6, 3).Value = 12752.5
